$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the event date in row 2 (19. 08. 2025 -> 20. 08. 2025)
$ws.Range("A2").Value = "20. 08. 2025"

# Capitalize event titles
$ws.Range("B6").Value = "Přednáška na náměstí"
$ws.Range("B9").Value = "Pečení"

# Update the saved cursor/selection position
$ws.Range("D21").Select()
